# Lista 1.docx - "Classes enum + Exercicio 14 Astah"
#
# The edit:
#  1. The empty paragraph right after the last table (which only held the
#     hidden "_GoBack" bookmark) becomes a plain empty paragraph - the
#     bookmark is relocated to the very end of the document.
#  2. One trailing space is trimmed from the end of paragraph "16- ...".
#  3. A new paragraph "19 - " (formatted like the surrounding answers) is
#     appended after paragraph "16- ...", carrying the relocated
#     "_GoBack" bookmark, followed by the pre-existing trailing empty
#     paragraph.

$d = $word.ActiveDocument

# --- Step 1: strip the _GoBack bookmark out of the paragraph right after
#     the table, leaving a plain empty paragraph in its place.
$bookmarkParaIndex = 154
$pOldBookmark = $d.Paragraphs.Item($bookmarkParaIndex)
$emptyParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pOldBookmark.Range.InsertXML($emptyParaXml) | Out-Null

# --- Step 2: remove one of the four trailing spaces after
#     "... pois há dependência de várias classes." Re-write just that
#     final run's text via InsertXML (a plain Range.Delete on a single
#     character ends up coalescing this run with the preceding ","
#     run in this host, which the reference diff does not do).
$answer16ParaIndex = 156
$pAnswer16 = $d.Paragraphs.Item($answer16ParaIndex)
$commaSearch = $pAnswer16.Range.Duplicate
$commaSearch.Find.Execute("Histórico,", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$lastRunStart = $commaSearch.End
$lastRunEnd = $pAnswer16.Range.End - 1
$lastRunRange = $d.Range($lastRunStart, $lastRunEnd)
$trimmedRunXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:lang w:val="pt"/></w:rPr><w:t xml:space="preserve"> pois há dependência de várias classes.   </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$lastRunRange.InsertXML($trimmedRunXml) | Out-Null

# --- Step 3: append a new paragraph "19 - " (same paragraph formatting as
#     the "15-"/"16-" answers) right after paragraph "16- ...", carrying
#     the "_GoBack" bookmark that used to sit right after the table.
$pAnswer16 = $d.Paragraphs.Item($answer16ParaIndex)
$pAnswer16.Range.InsertParagraphAfter() | Out-Null
$pNew = $d.Paragraphs.Item($answer16ParaIndex + 1)

$newParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="200" w:line="276" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:lang w:val="pt"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:lang w:val="pt"/></w:rPr><w:t xml:space="preserve">19 - </w:t></w:r><w:bookmarkStart w:id="6" w:name="_GoBack"/><w:bookmarkEnd w:id="6"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pNew.Range.InsertXML($newParaXml) | Out-Null
